$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "ОВУ" to "РХБЗ".
# All workbook-level definedName references (EmptyRows, RegisterDate, SoldierList,
# Заголовок, ИмяПодразделения, КВ, Подпись, Преподаватели) point at this sheet and
# get their "ОВУ!" qualifier updated to "РХБЗ!" automatically by Excel when the
# sheet is renamed.
$ws.Name = "РХБЗ"

# Switch the workbook calculation mode to Manual (xlCalculationManual = -4135).
$excel.Calculation = -4135
